# Simulated Wild Card round and logged it
# Update the "H" (Home) row target-depth stats on both the OFF and DEF
# sheets to reflect the additional Wild Card game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 ("H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 239
$wsOff.Range("C2").Value = 161
$wsOff.Range("D2").Value = 65
$wsOff.Range("E2").Value = 24

# --- DEF sheet: row 2 ("H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 181
$wsDef.Range("C2").Value = 126
$wsDef.Range("D2").Value = 41
$wsDef.Range("E2").Value = 22
$wsDef.Range("G2").Value = 2
